$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 353
$ws.Range("D6").Value = 282
$ws.Range("E6").Value = 71
$ws.Range("F6").Value = 62.80623608017817
$ws.Range("G6").Value = 20.11331444759207
$ws.Range("H6").Value = 79.88668555240793
